$d = $word.ActiveDocument

# Replace the standalone word "amount" with "number" in the sentence about
# backers. Using whole-word matching to avoid touching other occurrences.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "amount"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "number"
$find.Forward = $true
$find.Wrap = 1
$find.Format = $false
$find.MatchCase = $false
$find.MatchWholeWord = $true
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false
$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null
